$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D and E columns (Price / Volume change) as Text first so the
# replacement values keep their exact literal representation
# (e.g. "296.53", "-2.22%") instead of being auto-converted to numbers.
$ws.Range("D2:E47").NumberFormat = "@"

$ws.Range("D2").Value = "296.53"
$ws.Range("E2").Value = "-2.22%"
$ws.Range("D3").Value = "31.20"
$ws.Range("E3").Value = "-1.82%"
$ws.Range("D4").Value = "5.112"
$ws.Range("E4").Value = "-2.00%"
$ws.Range("D5").Value = "0.07335"
$ws.Range("E5").Value = "-0.31%"
$ws.Range("D6").Value = "7.714"
$ws.Range("E6").Value = "-1.68%"
$ws.Range("E7").Value = "12.19%"
$ws.Range("D8").Value = "3.733"
$ws.Range("E8").Value = "-0.02%"
$ws.Range("D9").Value = "0.9187"
$ws.Range("E9").Value = "1.18%"
$ws.Range("D10").Value = "0.1677"
$ws.Range("E10").Value = "-0.22%"
$ws.Range("D11").Value = "0.07024"
$ws.Range("E11").Value = "-6.37%"
$ws.Range("D12").Value = "0.08024"
$ws.Range("D13").Value = "0.02995"
$ws.Range("E13").Value = "1.40%"
$ws.Range("D14").Value = "0.09910"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").Value = "-0.27%"
$ws.Range("D16").Value = "0.006135"
$ws.Range("E16").Value = "-0.42%"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").Value = "-0.61%"
$ws.Range("D18").Value = "2.228"
$ws.Range("E18").Value = "-0.09%"
$ws.Range("D20").Value = "0.1331"
$ws.Range("E20").Value = "0.66%"
$ws.Range("D21").Value = "4.552"
$ws.Range("E21").Value = "0.47%"
$ws.Range("D22").Value = "0.04630"
$ws.Range("E22").Value = "2.15%"
$ws.Range("E23").Value = "-4.26%"
$ws.Range("E24").Value = "-0.12%"
$ws.Range("D25").Value = "0.004427"
$ws.Range("E25").Value = "0.14%"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").Value = "0.02%"
$ws.Range("D27").Value = "0.0001873"
$ws.Range("E27").Value = "7.74%"
$ws.Range("D39").Value = "0.01700"
$ws.Range("E39").Value = "1.25%"
$ws.Range("D40").Value = "0.04419"
$ws.Range("E40").Value = "-1.69%"
$ws.Range("D41").Value = "0.007203"
$ws.Range("E41").Value = "0.02%"
$ws.Range("D42").Value = "0.1328"
$ws.Range("E42").Value = "-1.21%"
$ws.Range("D43").Value = "0.002137"
$ws.Range("D44").Value = "0.01110"
$ws.Range("E44").Value = "-13.35%"
$ws.Range("D45").Value = "0.00005994"
$ws.Range("E45").Value = "-1.45%"
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").Value = "0.01021"
$ws.Range("E46").Value = "-21.32%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "1.868"
$ws.Range("E47").Value = "-1.27%"
